$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46060 -> 46061) for every data row (rows 2 through 171).
for ($r = 2; $r -le 171; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46060) {
        $cell.Value2 = 46061
    }
}
